$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy style from existing header cell (H1) to new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Data values for column I (I0) rows 2-7
$colI = @(5, 6, 9, 6, 7, 6)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $colI[$i]
}

# Data values for column J (IF) rows 2-7
$colJ = @(6, 6, 9, 6, 7, 6)
for ($i = 0; $i -lt $colJ.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $colJ[$i]
}
